$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) / Volume(1h) (E) updates for most rows. Values that look like a
# plain number (e.g. "580.63") are prefixed with a leading apostrophe so
# Excel stores them as text (matching the sheet's existing text-cell layout)
# instead of silently converting them to a numeric value.
$updates = @(
    @{ Row = 2; D = "67.663.22"; E = "  -0.62%  " },
    @{ Row = 3; D = "3.253.91"; E = "  -0.42%  " },
    @{ Row = 5; D = "'580.63"; E = "  -1.01%  " },
    @{ Row = 6; D = "'184.85"; E = "  +0.07%  " },
    @{ Row = 7; D = $null; E = "  -0.01%  " },
    @{ Row = 8; D = $null; E = "  +1.29%  " },
    @{ Row = 9; D = "3.250.56"; E = "  -0.46%  " },
    @{ Row = 10; D = $null; E = "  -2.50%  " },
    @{ Row = 11; D = "'6.56"; E = "  -2.38%  " },
    @{ Row = 12; D = $null; E = "  -0.98%  " },
    @{ Row = 13; D = "3.813.37"; E = "  -0.65%  " },
    @{ Row = 14; D = $null; E = "  -0.15%  " },
    @{ Row = 15; D = "'27.67"; E = "  -3.26%  " },
    @{ Row = 16; D = "67.656.65"; E = "  -0.61%  " },
    @{ Row = 17; D = $null; E = "  -1.40%  " },
    @{ Row = 18; D = "3.250.31"; E = "  -0.64%  " },
    @{ Row = 19; D = $null; E = "  -1.74%  " },
    @{ Row = 20; D = "'13.59"; E = "  -0.02%  " },
    @{ Row = 21; D = "'394.20"; E = "  +3.19%  " },
    @{ Row = 22; D = "'7.61"; E = "  -1.53%  " },
    @{ Row = 25; D = $null; E = "  +0.69%  " },
    @{ Row = 26; D = $null; E = "  -2.11%  " },
    @{ Row = 27; D = "'0.187"; E = "  -2.24%  " },
    @{ Row = 28; D = $null; E = "  -1.76%  " },
    @{ Row = 29; D = "'1.00"; E = "  +0.26%  " },
    @{ Row = 30; D = $null; E = "  -1.73%  " },
    @{ Row = 31; D = $null; E = "  -4.56%  " },
    @{ Row = 32; D = "'22.67"; E = "  -1.04%  " },
    @{ Row = 33; D = "'7.01"; E = "  -2.50%  " },
    @{ Row = 34; D = $null; E = "  -1.98%  " },
    @{ Row = 35; D = $null; E = "  +0.07%  " },
    @{ Row = 36; D = "'161.93"; E = "  -0.87%  " },
    @{ Row = 37; D = $null; E = "  -4.02%  " },
    @{ Row = 38; D = $null; E = "  +1.53%  " },
    @{ Row = 39; D = "'26.56"; E = "  -0.09%  " },
    @{ Row = 40; D = $null; E = "  -3.54%  " },
    @{ Row = 41; D = $null; E = "  -1.12%  " },
    @{ Row = 42; D = $null; E = "  -4.24%  " },
    @{ Row = 43; D = "'2.48"; E = "  -5.75%  " },
    @{ Row = 44; D = "'0.0689"; E = "  -0.12%  " },
    @{ Row = 45; D = "'40.67"; E = "  -1.25%  " },
    @{ Row = 46; D = "2.613.91"; E = "  -0.68%  " },
    @{ Row = 47; D = "'24.76"; E = "  -3.20%  " },
    @{ Row = 48; D = "'334.39"; E = "  -2.12%  " },
    @{ Row = 49; D = $null; E = "  -2.11%  " },
    @{ Row = 50; D = "'6.37"; E = "  +1.88%  " },
    @{ Row = 51; D = $null; E = "  -0.41%  " }
)

foreach ($u in $updates) {
    if ($u.D -ne $null) {
        $ws.Range("D" + $u.Row).Value = $u.D
    }
    $ws.Range("E" + $u.Row).Value = $u.E
}

# Rows 23/24 swap places: the coin formerly listed at row 23 (Dai) moves to
# row 24, and the coin formerly at row 24 (Litecoin) moves to row 23, each
# carrying its own refreshed price/volume figures.
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "'71.50"
$ws.Range("E23").Value = "  +0.17%  "

$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  +0.07%  "
